$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.908.73"
$ws.Range("E2").Value = "  -0.16%  "
Set-TextValue $ws.Range("D3") "1.874.61"
$ws.Range("E3").Value = "  -1.00%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.19%  "
Set-TextValue $ws.Range("D5") "0.7381"
$ws.Range("E5").Value = "  -4.74%  "
Set-TextValue $ws.Range("D6") "242.45"
$ws.Range("E6").Value = "  -0.53%  "
Set-TextValue $ws.Range("D7") "1.000"
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.3152"
$ws.Range("E8").Value = "  +0.80%  "
Set-TextValue $ws.Range("D9") "0.07161"
$ws.Range("E9").Value = "  -1.28%  "
Set-TextValue $ws.Range("D10") "24.65"
$ws.Range("E10").Value = "  -4.18%  "
Set-TextValue $ws.Range("D11") "0.08396"
$ws.Range("E11").Value = "  -3.59%  "
Set-TextValue $ws.Range("D12") "0.7500"
$ws.Range("E12").Value = "  -2.87%  "
Set-TextValue $ws.Range("D13") "5.420"
$ws.Range("E13").Value = "  +0.09%  "
Set-TextValue $ws.Range("D14") "1.840.45"
$ws.Range("E14").Value = "  -7.20%  "
Set-TextValue $ws.Range("D15") "92.51"
$ws.Range("E15").Value = "  -2.05%  "
Set-TextValue $ws.Range("D16") "29.882.77"
$ws.Range("E16").Value = "  -1.19%  "
Set-TextValue $ws.Range("D17") "6.100"
$ws.Range("E17").Value = "  -1.87%  "
Set-TextValue $ws.Range("D18") "13.58"
$ws.Range("E18").Value = "  -2.47%  "
Set-TextValue $ws.Range("D19") "243.01"
$ws.Range("E19").Value = "  -1.03%  "
Set-TextValue $ws.Range("D20") "0.000007809"
$ws.Range("E20").Value = "  -1.05%  "
Set-TextValue $ws.Range("D21") "0.9993"
$ws.Range("E21").Value = "  -0.04%  "
Set-TextValue $ws.Range("D22") "2.118.81"
$ws.Range("E22").Value = "  -11.10%  "
Set-TextValue $ws.Range("D23") "7.985"
$ws.Range("E23").Value = "  -2.24%  "
Set-TextValue $ws.Range("D24") "1.002"
$ws.Range("E24").Value = "  +0.20%  "
Set-TextValue $ws.Range("D25") "0.1548"
$ws.Range("E25").Value = "  -2.99%  "
Set-TextValue $ws.Range("D26") "9.292"
$ws.Range("E26").Value = "  -2.53%  "
Set-TextValue $ws.Range("D27") "165.10"
$ws.Range("E27").Value = "  +1.52%  "
Set-TextValue $ws.Range("D28") "18.59"
$ws.Range("E28").Value = "  -1.30%  "
Set-TextValue $ws.Range("D29") "2.035"
$ws.Range("E29").Value = "  -0.67%  "
Set-TextValue $ws.Range("D30") "1.491"
$ws.Range("E30").Value = "  +4.19%  "
Set-TextValue $ws.Range("D31") "4.597"
$ws.Range("E31").Value = "  +1.66%  "
Set-TextValue $ws.Range("D32") "1.533"
$ws.Range("E32").Value = "  -0.79%  "
Set-TextValue $ws.Range("D33") "4.254"
$ws.Range("E33").Value = "  +3.11%  "
Set-TextValue $ws.Range("D34") "0.05322"
$ws.Range("E34").Value = "  -2.20%  "
Set-TextValue $ws.Range("D35") "1.235"
$ws.Range("E35").Value = "  -1.22%  "
Set-TextValue $ws.Range("D36") "0.7536"
$ws.Range("E36").Value = "  +0.20%  "
Set-TextValue $ws.Range("D37") "0.9956"
$ws.Range("E37").Value = "  -0.90%  "
Set-TextValue $ws.Range("D38") "2.691"
$ws.Range("E38").Value = "  +0.19%  "
Set-TextValue $ws.Range("D39") "0.01950"
$ws.Range("E39").Value = "  -1.72%  "
Set-TextValue $ws.Range("D40") "2.755"
$ws.Range("E40").Value = "  -1.05%  "
Set-TextValue $ws.Range("D41") "0.4504"
$ws.Range("E41").Value = "  -0.26%  "
Set-TextValue $ws.Range("D42") "1.111.43"
$ws.Range("E42").Value = "  +1.44%  "
Set-TextValue $ws.Range("D43") "6.043"
$ws.Range("E43").Value = "  -0.73%  "
Set-TextValue $ws.Range("D44") "72.13"
$ws.Range("E44").Value = "  -1.81%  "
Set-TextValue $ws.Range("D45") "0.8579"
$ws.Range("E45").Value = "  +0.47%  "
Set-TextValue $ws.Range("D46") "1.001"
$ws.Range("E46").Value = "  +0.16%  "
Set-TextValue $ws.Range("D47") "103.08"
$ws.Range("E47").Value = "  +0.02%  "
Set-TextValue $ws.Range("D48") "7.646"
$ws.Range("E48").Value = "  +0.29%  "
Set-TextValue $ws.Range("D49") "3.095"
$ws.Range("E49").Value = "  +3.16%  "
Set-TextValue $ws.Range("D50") "1.838"
$ws.Range("E50").Value = "  -2.67%  "
Set-TextValue $ws.Range("D51") "2.014.52"
$ws.Range("E51").Value = "  -9.87%  "
